# Update simulation result values across the four result sheets.
# (Terminology refactor: 'manufactured death' -> 'structural barriers' model
#  re-run produced updated stochastic simulation outputs.)

$wb = $excel.ActiveWorkbook

# --- Sheet: National Forecast ---
$ws1 = $wb.Worksheets.Item("National Forecast")
$ws1.Range("B2").Value = 2000
$ws1.Range("B3").Value = 1943
$ws1.Range("B4").Value = 0.9715
$ws1.Range("B5").Value = 3
$ws1.Range("B6").Value = 3.792074112197632
$ws1.Range("B7").Value = 3.709746932931708
$ws1.Range("B10").Value = 6
$ws1.Range("B11").Value = 9
$ws1.Range("B12").Value = 0.7380339680905815
$ws1.Range("B13").Value = 0.9274318064848173

# --- Sheet: Regional Comparison ---
$ws2 = $wb.Worksheets.Item("Regional Comparison")
$ws2.Range("C2").Value = 976
$ws2.Range("D2").Value = 0.976
$ws2.Range("F2").Value = 0.7838114754098361

$ws2.Range("C3").Value = 996
$ws2.Range("D3").Value = 0.996
$ws2.Range("F3").Value = 0.8634538152610441

$ws2.Range("C4").Value = 988
$ws2.Range("D4").Value = 0.988
$ws2.Range("F4").Value = 0.7834008097165992

$ws2.Range("C5").Value = 971
$ws2.Range("D5").Value = 0.971
$ws2.Range("F5").Value = 0.7250257466529351

# --- Sheet: Scenario Comparison ---
$ws3 = $wb.Worksheets.Item("Scenario Comparison")
$ws3.Range("B2").Value = 0.6953846153846154
$ws3.Range("B3").Value = 0.6764091858037579
$ws3.Range("B4").Value = 0.700312174817898
$ws3.Range("C4").Value = 3
$ws3.Range("B5").Value = 0.6189451022604952
$ws3.Range("B6").Value = 0.6434782608695652
$ws3.Range("B7").Value = 0.5082382762991128

# --- Sheet: PSA Results ---
$ws4 = $wb.Worksheets.Item("PSA Results")
$ws4.Range("B2").Value = 0.7304657182353299
$ws4.Range("C2").Value = 0.6353124999999999
$ws4.Range("D2").Value = 0.82

$ws4.Range("B3").Value = 0.9252023740057934
$ws4.Range("C3").Value = 0.8699343434343434
$ws4.Range("D3").Value = 0.9696969696969697

$ws4.Range("B4").Value = 2.782
$ws4.Range("D4").Value = 4
